# Update quantity column (C) - plain numeric values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 67
$ws.Range("C9").Value = 65
$ws.Range("C10").Value = 58
$ws.Range("C11").Value = 9
$ws.Range("C12").Value = 74
$ws.Range("C13").Value = 80
$ws.Range("C14").Value = 67
$ws.Range("C15").Value = 54
$ws.Range("C16").Value = 79
$ws.Range("C17").Value = 83

# Update amount columns (G/H) - these are stored as text (e.g. "16640.00"),
# so force text number format before assigning so Excel does not
# auto-convert the numeric-looking string back into a number.
$amountCells = "G9","G10","G11","G13","G14","G19","H19","G21","H21"
foreach ($cellRef in $amountCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("G9").Value = "16640.00"
$ws.Range("G10").Value = "27376.00"
$ws.Range("G11").Value = "5958.00"
$ws.Range("G13").Value = "10880.00"
$ws.Range("G14").Value = "1541.00"
$ws.Range("G19").Value = "62395.00"
$ws.Range("H19").Value = "62395.00"
$ws.Range("G21").Value = "62395.00"
$ws.Range("H21").Value = "62395.00"
